$p = $ppt.ActivePresentation

# Add a new blank slide (layout 7 = "Blank") as slide 2
$p.Slides.Add(2, 12)
